$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.466.34"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "3.391.12"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'575.50"
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("D6").Value = "'140.55"
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("D9").Value = "'7.67"
$ws.Range("E9").Value = "  +1.79%  "
$ws.Range("E10").Value = "  -1.20%  "
$ws.Range("D11").Value = "'0.388"
$ws.Range("E11").Value = "  -1.73%  "
$ws.Range("D12").Value = "3.971.20"
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("D13").Value = "'0.126"
$ws.Range("E13").Value = "  +0.81%  "
$ws.Range("D14").Value = "'28.33"
$ws.Range("E14").Value = "  -0.38%  "
$ws.Range("D15").Value = "3.394.18"
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("E16").Value = "  -0.76%  "
$ws.Range("D17").Value = "61.464.34"
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("D18").Value = "'6.16"
$ws.Range("E18").Value = "  -1.01%  "
$ws.Range("D19").Value = "'13.67"
$ws.Range("E19").Value = "  -2.25%  "
$ws.Range("D20").Value = "'8.98"
$ws.Range("E20").Value = "  -1.05%  "
$ws.Range("D21").Value = "'389.45"
$ws.Range("E21").Value = "  +0.56%  "
$ws.Range("D22").Value = "'75.22"
$ws.Range("E22").Value = "  +1.42%  "
$ws.Range("E23").Value = "  -0.79%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("E25").Value = "  -4.43%  "
$ws.Range("E26").Value = "  +6.63%  "
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("D28").Value = "'7.30"
$ws.Range("E28").Value = "  -1.70%  "
$ws.Range("D29").Value = "'8.02"
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("E30").Value = "  -0.38%  "
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").Value = "'1.38"
$ws.Range("E32").Value = "  -3.44%  "
$ws.Range("D33").Value = "'23.40"
$ws.Range("E33").Value = "  -0.99%  "
$ws.Range("D34").Value = "'6.92"
$ws.Range("E34").Value = "  -1.58%  "
$ws.Range("D35").Value = "'167.50"
$ws.Range("E35").Value = "  +0.37%  "
$ws.Range("E36").Value = "  +0.52%  "
$ws.Range("D37").Value = "3.426.60"
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").Value = "'1.48"
$ws.Range("E38").Value = "  -0.89%  "
$ws.Range("D39").Value = "'0.0768"
$ws.Range("E39").Value = "  -1.36%  "
$ws.Range("D40").Value = "'25.79"
$ws.Range("E40").Value = "  -10.27%  "
$ws.Range("E41").Value = "  -0.46%  "
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("E43").Value = "  -0.68%  "
$ws.Range("E44").Value = "  -0.25%  "
$ws.Range("D45").Value = "2.449.59"
$ws.Range("E45").Value = "  -2.03%  "
$ws.Range("D46").Value = "'22.79"
$ws.Range("E46").Value = "  -2.26%  "
$ws.Range("E47").Value = "  -1.89%  "
$ws.Range("D48").Value = "'0.999"
$ws.Range("E48").Value = "  -0.10%  "
$ws.Range("E49").Value = "  -3.05%  "
$ws.Range("D50").Value = "'2.04"
$ws.Range("E50").Value = "  -2.23%  "
$ws.Range("E51").Value = "  -1.86%  "
